$d = $word.ActiveDocument
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---- Helper: find paragraph index by exact (trimmed) text ----
function Find-ParaIndex($doc, $text) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $pt = $doc.Paragraphs($i).Range.Text.TrimEnd("`r")
        if ($pt -eq $text) {
            return $i
        }
    }
    return -1
}

# =====================================================================
# 1) "these kind of problems" -> split into two runs reading
#    "... for " / "these kinds of problems"
# =====================================================================
$i8 = Find-ParaIndex $d "The integration of machine learning (ML) into medical imaging offers a promising solution for these kind of problems. Deep learning-based machine learning algorithms, in particular, are capable of processing enormous volumes of imaging data and spotting intricate patterns that could be invisible to human observers. According to Litjens et al. (2017), these algorithms may decrease variability, increase diagnostic accuracy, and support clinical decision making."
$para8 = $d.Paragraphs($i8)
$rngFull = $para8.Range
$fullText = $rngFull.Text
$searchStr = "these kind of problems"
$idx = $fullText.IndexOf($searchStr)
$subStart = $rngFull.Start + $idx
$subEnd = $subStart + $searchStr.Length
$sub = $d.Range($subStart, $subEnd)
$sub.Text = "these kinds of problems"
# Force a run split at the same formatting (no visual change) so the
# fixed phrase ends up as its own run, matching the target markup.
$subEnd2 = $subStart + ("these kinds of problems").Length
$sub2 = $d.Range($subStart, $subEnd2)
$sub2.Font.Size = 13
$sub2.Font.Size = 12

# =====================================================================
# 2) Insert five new paragraphs right after that paragraph (before the
#    "Problem statement" heading)
# =====================================================================
$para8b = $d.Paragraphs($i8)
$endPt = $d.Range($para8b.Range.End - 1, $para8b.Range.End - 1)

$rPr = "<w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>"

$newParasXml = @"
<w:p $w>
  <w:pPr>$rPr</w:pPr>
  <w:r>$rPr<w:t>Machine learning comprises a wide range of approaches that are divided into two categories supervised and unsupervised learning. In supervised learning, models are trained using labelled data to predict or categorise the data. On the other hand, unsupervised learning entails identifying structures or hidden patterns in unlabelled data. The medical industry can benefit greatly from both forms of learning, particularly in the areas of cancer detection and diagnosis.</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>$rPr</w:pPr>
  <w:r>$rPr<w:t xml:space="preserve">In the context of prostate cancer, machine learning models have been developed for various tasks such as tumour detection, Gleason grade prediction, and treatment response monitoring. Recurrent neural networks (RNN), convolutional neural network (CNN), and </w:t></w:r>
  <w:r>$rPr<w:lastRenderedPageBreak/><w:t>support vector machines (SVM) are notable machine learning approaches. To varying extents these models have improved prostate cancer diagnosis efficiency and accuracy</w:t></w:r>
  <w:r>$rPr<w:t>.</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>$rPr</w:pPr>
  <w:r>$rPr<w:t>For instance, since CNNs can learn spatial hierarchies from input images, they are particularly useful for image analysis tasks like tumour segmentation and classification (Pellicer-Valero et al., 2022). Similar to this, RNNs which are capable at preprocessing sequential data have been used to forecast treatment results by taking into account the history and advancements of patients (Mirsamadi et al., 2017).</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>$rPr</w:pPr>
  <w:r>$rPr<w:t xml:space="preserve">Even with the improvements, there are still a number of difficulties in applying ML to clinical practice. These include obtaining clinical validation to verify the </w:t></w:r>
  <w:r>$rPr<w:t>model&#8217;s</w:t></w:r>
  <w:r>$rPr<w:t xml:space="preserve"> efficiency in real </w:t></w:r>
  <w:r>$rPr<w:t>world</w:t></w:r>
  <w:r>$rPr<w:t xml:space="preserve"> </w:t></w:r>
  <w:r>$rPr<w:t>scenarios, guaranteeing the interpretability and transparency of ML decisions, and training robust models on big annotated datasets. To tackle these obstacles, data scientists, physicians, and regulatory agencies must continue their research and work together.</w:t></w:r>
</w:p>
<w:p $w>
  <w:pPr>$rPr</w:pPr>
  <w:r>$rPr<w:t>The development of sophisticated imaging technologies and machine learning presents a possible alternative to the poor performance of existing methods for the identification of prostate cancer. Utilising these technologies can lead to better patient outcomes by enhancing diagnostic accuracy and personalising treatment</w:t></w:r>
  <w:r>$rPr<w:t xml:space="preserve"> approaches. To overcome current obstacles and fully realise the potential of these state of the art instruments in the treatment of prostate cancer, more research and innovation in this area are imperative.</w:t></w:r>
</w:p>
"@

$endPt.InsertXML($newParasXml)

# =====================================================================
# 3) Move the lastRenderedPageBreak markers caused by the extra page:
#    - remove it from "Chapter 2 Literature review"
#    - add it before "2.4 Related work"
#    - remove it from "3.5 Model Development"
#    - add it before "4.3 Comparison with the existing methods"
# =====================================================================
$iCh2 = Find-ParaIndex $d "Chapter 2 Literature review"
$pCh2 = $d.Paragraphs($iCh2)
$pCh2.Range.InsertXML("<w:p $w><w:pPr><w:pStyle w:val='Heading1'/></w:pPr><w:r><w:t>Chapter 2 Literature review</w:t></w:r></w:p>")

$i24 = Find-ParaIndex $d "2.4 Related work"
$p24 = $d.Paragraphs($i24)
$p24.Range.InsertXML("<w:p $w><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>2.4 Related work</w:t></w:r></w:p>")

$i35 = Find-ParaIndex $d "3.5 Model Development"
$p35 = $d.Paragraphs($i35)
$p35.Range.InsertXML("<w:p $w><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:t>3.5 Model Development</w:t></w:r></w:p>")

$i43 = Find-ParaIndex $d "4.3 Comparison with the existing methods"
$p43 = $d.Paragraphs($i43)
$p43.Range.InsertXML("<w:p $w><w:pPr><w:pStyle w:val='Heading2'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>4.3 Comparison with the existing methods</w:t></w:r></w:p>")

Write-Output "done"
